# "Drache Sound + Weltraumtheme added"
# Adds the "Weltraum Theme (Dark)" entry into the previously empty row 23
# of the sound planning table (left block B:F and its mirrored right
# block H:L), and updates the sheet's zoom/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Left table (Soundbezeichnung / Zeitpunkt / Aufgabe / Details / Priorität)
$ws.Range("B23").Value = "Weltraum Theme (Dark)"
$ws.Range("D23").Value = "fertig zum Einfügen in Game"
$ws.Range("E23").Value = "Loop-fähig"
$ws.Range("F23").Value = 10

# Mirrored right table only gets the name + priority filled in
$ws.Range("H23").Value = "Weltraum Theme (Dark)"
$ws.Range("L23").Value = 10

# Sheet view was zoomed out and the selection moved onto the new row
$ws.Select()
$excel.ActiveWindow.Zoom = 70
$ws.Range("E23").Select()
